# Remove the "COMPLEX SHEET" tab from the faculty enrollment import
# template, keeping only the simple sheet (renamed to "FACULTY USERS").

$wb = $excel.ActiveWorkbook

# Rename the remaining (simple) sheet.
$ws = $wb.Worksheets.Item(1)
$ws.Name = "FACULTY USERS"

# Delete the complex sheet tab - suppress the "delete sheet" confirmation.
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("COMPLEX SHEET").Delete()

# Move the active selection to C6, matching the saved view state.
$ws.Range("C6").Select()
